$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.481.03"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "  +7.04%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.624.83"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = "  +7.00%  "

$ws.Range("E4").Value = "  +0.29%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.98"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = "  +5.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "192.02"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = "  +8.98%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.648"
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = "  +2.79%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.600.57"
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = "  +6.67%  "

$ws.Range("E9").Value = "  +0.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.180"
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = "  +4.27%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.661"
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = "  +4.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "57.66"
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = "  +7.52%  "

$ws.Range("E13").Value = "  +6.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.73"
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = "  +5.53%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.212.72"
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = "  +7.54%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.632.66"
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = "  +7.61%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.32"
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = "  +5.93%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.465.83"
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = "  +7.21%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.57"
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = "  +5.89%  "

$ws.Range("E20").Value = "  +0.97%  "

$ws.Range("E21").Value = "  +4.82%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "493.69"
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = "  +5.42%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.51"
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = "  +12.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.68"
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = "  +17.10%  "

$ws.Range("E25").Value = "  +8.48%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "90.48"
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = "  +0.81%  "

$ws.Range("E27").Value = "  +6.38%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.19"
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = "  +5.41%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.35"
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = "  +7.21%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.35"
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = "  +3.82%  "

$ws.Range("E31").Value = "  +15.24%  "

$ws.Range("E32").Value = "  +6.80%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "613.89"
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = "  +6.36%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "65.32"
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = "  +5.19%  "

$ws.Range("E35").Value = "  +7.04%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0829"
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = "  +12.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.148"
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = "  +4.47%  "

$ws.Range("E38").Value = "  -0.05%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.85"
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = "  +5.24%  "

$ws.Range("E40").Value = "  +6.69%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.68"
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = "  +1.72%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.357.06"
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = "  +8.34%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.05"
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = "  +7.57%  "

$ws.Range("E44").Value = "  +6.82%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.67"
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = "  +9.21%  "

$ws.Range("E46").Value = "  +6.44%  "

$ws.Range("E47").Value = "  +2.65%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.74"
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = "  +10.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.07"
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = "  +6.74%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.30"
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = "  +4.72%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = "  +0.32%  "
